$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing A:H data shifts to B:I.
$ws.Columns("A").Insert()

# New "Id" header in A2, matching the plain-bordered (non-bold) header style.
$ws.Range("A2").Value = "Id"
$ws.Range("A2").Borders.LineStyle = 1

# Sequential id values for the data rows (now rows 3-12).
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(3 + $i, 1).Value = $i + 1
}

# Leftover formatting on C13 (underlined, empty cell) from the edit session.
$ws.Range("C13").Font.Underline = 2

# Match the final selection left in the saved workbook.
[void]$ws.Range("C13").Select()

Write-Output "ok"
